# Weekly price update: a new record for the week of 2022-05-25 is inserted
# at the top of the "Femacal de La Calera - Alcachofa" block (row 279),
# pushing all existing rows in that block down by one (279-297 -> 280-298).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 279; this shifts rows 279:297 down to 280:298
# and carries the D-column date style (s="2") onto the new blank row.
$ws.Rows.Item(279).Insert()

# Populate the new row 279 with the latest weekly record.
$ws.Cells.Item(279, 1).Value = 3
$ws.Cells.Item(279, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(279, 3).Value = "Coquimbo"
$ws.Cells.Item(279, 4).Value = 44706
$ws.Cells.Item(279, 5).Value = 5
$ws.Cells.Item(279, 6).Value = 100112013
$ws.Cells.Item(279, 7).Value = "Alcachofa"
$ws.Cells.Item(279, 8).Value = "Argentina(o)"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 110
$ws.Cells.Item(279, 11).Value = 16500
$ws.Cells.Item(279, 12).Value = 17000
$ws.Cells.Item(279, 13).Value = 16727
$ws.Cells.Item(279, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(279, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(279, 16).Value = 335
$ws.Cells.Item(279, 17).Value = 50
$ws.Cells.Item(279, 18).Value = "Hortaliza"
